# Natmi following Dr Hou advice
#
# Bmp3 -> Acvr2b ligand-receptor table: the sending/target cluster labels
# are updated (an extra cell population, "sCs", is split into "sCs" and a
# new "ECs" cluster) and every LR-pair statistic is recomputed. The single
# "FAPs -> FAPs / FAPs -> sCs" pair (rows 2-3) becomes six sender/target
# combinations spanning FAPs and the newly distinguished sCs (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row (row 1 is the header, already correct).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the LR-pair statistics.
$rows = @(
    @{ Row=2; Values=@("FAPs","Bmp3","Acvr2b","ECs",
            3, 1, 3.131909666666667, 9.395728999999999,
            0.9475081863095236, 0.9475081863095236, 2, 0.6666666666666666,
            1.539665666666667, 4.618997, 0.3572088291809875, 0.3572088291809875,
            4.822093784868112, 43.398844063813, 0.3384582898710259, 0.3384582898710259) }
    @{ Row=3; Values=@("FAPs","Bmp3","Acvr2b","FAPs",
            3, 1, 3.131909666666667, 9.395728999999999,
            0.9475081863095236, 0.9475081863095236, 3, 1,
            1.452872333333333, 4.358617, 0.3370724153789985, 0.3370724153789985,
            4.550264905199222, 40.95238414679299, 0.3193788729507253, 0.3193788729507253) }
    @{ Row=4; Values=@("FAPs","Bmp3","Acvr2b","sCs",
            3, 1, 3.131909666666667, 9.395728999999999,
            0.9475081863095236, 0.9475081863095236, 3, 1,
            1.317729666666667, 3.953189, 0.3057187554400141, 0.3057187554400141,
            4.127010281086778, 37.143092529781, 0.2896710234877726, 0.2896710234877726) }
    @{ Row=5; Values=@("sCs","Bmp3","Acvr2b","ECs",
            1, 0.3333333333333333, 0.1735073333333333, 0.520522,
            0.05249181369047638, 0.05249181369047638, 2, 0.6666666666666666,
            1.539665666666667, 4.618997, 0.3572088291809875, 0.3572088291809875,
            0.2671432840482223, 2.404289556434, 0.0187505393099616, 0.0187505393099616) }
    @{ Row=6; Values=@("sCs","Bmp3","Acvr2b","FAPs",
            1, 0.3333333333333333, 0.1735073333333333, 0.520522,
            0.05249181369047638, 0.05249181369047638, 3, 1,
            1.452872333333333, 4.358617, 0.3370724153789985, 0.3370724153789985,
            0.2520840042304445, 2.268756038074, 0.01769354242827325, 0.01769354242827326) }
    @{ Row=7; Values=@("sCs","Bmp3","Acvr2b","sCs",
            1, 0.3333333333333333, 0.1735073333333333, 0.520522,
            0.05249181369047638, 0.05249181369047638, 3, 1,
            1.317729666666667, 3.953189, 0.3057187554400141, 0.3057187554400141,
            0.2286357605175556, 2.057721844658, 0.01604773195224153, 0.01604773195224153) }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
